# Exp3Groups.xlsx update
# - Add "ID" header in Sheet1!A1
# - Update several Perf values in column B
# - Make Sheet1 the active/selected sheet (instead of Terms) and move the
#   selection/view down to around row 99

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add new "ID" header cell
$ws1.Range("A1").Value = "ID"

# Corrected / updated performance values
$ws1.Range("B9").Value  = 3.77
$ws1.Range("B16").Value = 4.0999999999999996
$ws1.Range("B23").Value = 3.84
$ws1.Range("B25").Value = 3.84
$ws1.Range("B32").Value = 4.0999999999999996
$ws1.Range("B43").Value = 2.3199999999999998
$ws1.Range("B51").Value = 2.38
$ws1.Range("B61").Value = 2.25
$ws1.Range("B83").Value = 4.45
$ws1.Range("B85").Value = 4.51
$ws1.Range("B95").Value = 4.47
$ws1.Range("B108").Value = 4.45

# Switch active sheet/selection: Sheet1 becomes the selected tab,
# with the view scrolled down and B99 selected; Terms loses its
# tabSelected flag (its own A10 selection is left untouched).
$ws1.Activate() | Out-Null
$ws1.Range("B99").Select() | Out-Null
